$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FIFO")

# Insert a new row above the old row 3 (FIFO_W_01) to hold the new
# FIFO_INIT_02 requirement, pushing everything else down by one.
$ws.Rows("3:3").Insert()

# Row 2 (FIFO_INIT_01): reworded requirement body.
$ws.Range("D2").Value = "Module shall set FIFO_DEPTH to default 2 (optionally 2,4,8,16,32,64) upon initialization."

# New row 3: FIFO_INIT_02 requirement.
$ws.Range("C3").Value = "FIFO_INIT_02"
$ws.Range("D3").Value = "Module shall include internal memory buffer mem_array to hold FIFO data"

# Re-merge the "Initialization" function cell over its two rows now that a
# new row was inserted beneath it.
$ws.Range("B2:B3").Merge()

# Reworded requirement bodies for the rows that shifted down by one.
$ws.Range("D5").Value = "Module shall write 128 bits of o_fifo_w_data to mem_array[w_ptr] when i_w_en is logic high."
$ws.Range("D6").Value = "Module shall set o_full logic high when w_ptr is equal to 1 less than the buffer length and r_ptr is 0."
$ws.Range("D8").Value = "Module shall set 128 bits of o_fifo_r_data to value of mem_array[r_ptr] when i_r_en is logic high."
$ws.Range("D9").Value = "Module shall set o_empty logic high when r_ptr is equal to 1 less than the buffer length and w_ptr is 0."

# Restore selection to the cell the author left active.
$ws.Range("D5").Select()
